# Commit: Updated getWatershed to use IEc file, updated tablesOut and
# tablesToExcel to reflect the change from Upper and Lower Mortendad.
#
# This collapses the "Upper Mortendad" / "Lower Mortendad" watershed split
# into a single "Mortendad" watershed, and (since the well R-5 SC's
# watershed is no longer resolved the same way under the new IEc-file-based
# getWatershed) drops the R-5 SC row from both the mapping table and the
# exhibit table.

$wb = $excel.ActiveWorkbook

# --- Sheet "Regional for Mapping" -----------------------------------------
$wsMap = $wb.Worksheets.Item("Regional for Mapping")

# P3 ("Lower Mortendad") -> "Mortendad" (Upper/Lower split removed)
$wsMap.Cells.Item(3, 16).Value2 = "Mortendad"

# Row 5 is the R-5 SC well; it no longer belongs in this table, so remove
# the entire row (rows below shift up, e.g. old row 6 "R-8 OB" becomes the
# new row 5, etc.)
$wsMap.Rows.Item(5).Delete()

# --- Sheet "Regional Exhibit" ---------------------------------------------
$wsExhibit = $wb.Worksheets.Item("Regional Exhibit")

# Row 7 section header "Mortendad Canyon" -> "Los Alamos and Pajarito Canyons"
$wsExhibit.Cells.Item(7, 1).Value2 = "Los Alamos and Pajarito Canyons"

# Row 9 is the R-5 SC well entry; remove it so later rows shift up.
$wsExhibit.Rows.Item(9).Delete()
